{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" boilerplate\n// paragraphs (and the blank paragraph that separated them from the\n// requirements line) that the site generator stopped emitting.\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = body.paragraphs.items;\n\n// Locate the anchor paragraph (\"LOQ4205: ... (Requisito fraco)\") that\n// precedes the block being removed.\nlet anchorIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"LOQ4205\") !== -1) {\n    anchorIndex = i;\n    break;\n  }\n}\n\nif (anchorIndex === -1) {\n  throw new Error(\"Could not locate the 'LOQ4205' requirements paragraph.\");\n}\n\n// The three paragraphs right after the anchor are:\n//   1) an empty spacer paragraph\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n//   3) \"\u00a9 2020 . Contact: ... Creative Commons Attribution\"\n// Delete them (walking via getNext so we don't depend on fixed indices).\nlet toDelete = [];\nlet cursor = items[anchorIndex].getNext();\nfor (let n = 0; n < 3; n++) {\n  cursor.load(\"text\");\n  toDelete.push(cursor);\n  cursor = cursor.getNext();\n}\nawait context.sync();\n\nconst expectedTexts = [\n  \"\",\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n];\n\nfor (let i = 0; i < toDelete.length; i++) {\n  if (toDelete[i].text !== expectedTexts[i]) {\n    throw new Error(\"Unexpected paragraph content while deleting boilerplate: \" + toDelete[i].text);\n  }\n}\n\nfor (const para of toDelete) {\n  para.delete();\n}\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"(c) 2020 ...\" boilerplate\n# paragraphs (and the blank spacer paragraph that separated them from the\n# requirements line) that the site generator stopped emitting.\n\n$d = $word.ActiveDocument\n\n# Locate the anchor paragraph (\"LOQ4205: ... (Requisito fraco)\") that\n# precedes the block being removed.\n$count = $d.Paragraphs.Count\n$anchorIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*LOQ4205*\") {\n        $anchorIndex = $i\n        break\n    }\n}\nif ($anchorIndex -eq -1) {\n    throw \"Could not locate the 'LOQ4205' requirements paragraph.\"\n}\n\n# The three paragraphs right after the anchor are:\n#   1) an empty spacer paragraph\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n#   3) \"(c) 2020 . Contact: ... Creative Commons Attribution\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n$expectedTexts = @(\n    \"\",\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    $copyrightText\n)\n\nfor ($n = 0; $n -lt 3; $n++) {\n    # Each deletion shifts everything after it up by one, so the\n    # paragraph we want is always the one right after the anchor.\n    $target = $d.Paragraphs.Item($anchorIndex + 1)\n    $text = $target.Range.Text.TrimEnd([char]13, [char]7)\n    if ($text -ne $expectedTexts[$n]) {\n        throw \"Unexpected paragraph content while deleting boilerplate: $text\"\n    }\n    $target.Range.Delete()\n}\n"}
